$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add the two new rows of data (URL / Note) at the bottom of the table.
$ws.Range("A8").Value = "http://www2.compute.dtu.dk/~pcha/AIRtools/AIRtoolsManual.pdf"
$ws.Range("B8").Value = "Tool box for matlab"
$ws.Range("A9").Value = "https://tomroelandts.com/articles/do-not-ignore-the-astra-toolbox"

# Move the active selection to A8, matching the saved workbook view state.
$ws.Range("A8").Select()
